$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 62

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "hassaan"
$ws.Range("C5").Value = "Karachi"
$ws.Range("D5").Value = "Gold"
$ws.Range("E5").Value = 52
$ws.Range("F5").Value = 5.4
$ws.Range("G5").Value = "2025-03-17"
$ws.Range("H5").Value = "Paid"
$ws.Range("K5").Value = "2025-02-25"
$ws.Range("M5").Value = "Male"
$ws.Range("R5").Value = "strong"
$ws.Range("S5").Value = "u"
